$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.889.91'
$ws.Range("D3").Value = '3.868.36'
$ws.Range("E3").Value = '  +2.96%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.11'
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.63'
$ws.Range("E6").Value = '  -2.79%  '
$ws.Range("D7").Value = '3.867.81'
$ws.Range("E7").Value = '  +2.96%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("E11").Value = '  -1.56%  '
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.01'
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").Value = '4.514.45'
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("D16").Value = '3.867.39'
$ws.Range("E16").Value = '  +3.21%  '
$ws.Range("D17").Value = '69.059.81'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.56'
$ws.Range("E18").Value = '  +2.59%  '
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.35'
$ws.Range("E20").Value = '  +2.56%  '
$ws.Range("E21").Value = '  -1.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '485.09'
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.720'
$ws.Range("E23").Value = '  -1.42%  '
$ws.Range("E24").Value = '  +7.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.01'
$ws.Range("E25").Value = '  -1.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.25'
$ws.Range("E26").Value = '  -2.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.12'
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.98'
$ws.Range("E29").Value = '  -1.21%  '
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.93'
$ws.Range("E31").Value = '  -2.84%  '
$ws.Range("D32").Value = '4.015.90'
$ws.Range("E32").Value = '  +2.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.35'
$ws.Range("E33").Value = '  +2.30%  '
$ws.Range("E34").Value = '  -4.30%  '
$ws.Range("D35").Value = '3.814.19'
$ws.Range("E35").Value = '  +3.36%  '
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("E38").Value = '  +1.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.90'
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E41").Value = '  -2.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '437.44'
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.99'
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.41'
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '143.46'
$ws.Range("E48").Value = '  +1.56%  '
$ws.Range("D49").Value = '2.841.99'
$ws.Range("E49").Value = '  +1.69%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.21'
$ws.Range("E50").Value = '  +13.98%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0359'
$ws.Range("E51").Value = '  +1.51%  '
